$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

$ws.Range("D2").Value = "62.281.70"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").Value = "3.429.83"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("E4").Value = "  -0.12%  "
Set-TextValue $ws "D5" "579.17"
$ws.Range("E5").Value = "  -0.35%  "
Set-TextValue $ws "D6" "152.65"
$ws.Range("E6").Value = "  +3.55%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +1.20%  "
Set-TextValue $ws "D9" "8.05"
$ws.Range("E9").Value = "  +4.25%  "
Set-TextValue $ws "D10" "0.125"
$ws.Range("E10").Value = "  -0.51%  "
Set-TextValue $ws "D11" "0.417"
$ws.Range("E11").Value = "  +2.72%  "
$ws.Range("D12").Value = "4.022.42"
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("E13").Value = "  +0.27%  "
Set-TextValue $ws "D14" "28.74"
$ws.Range("E14").Value = "  -3.51%  "
$ws.Range("D15").Value = "3.433.52"
$ws.Range("E15").Value = "  -1.36%  "
Set-TextValue $ws "D16" "0.0000172"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "62.266.08"
$ws.Range("E17").Value = "  -2.03%  "
Set-TextValue $ws "D18" "6.50"
$ws.Range("E18").Value = "  +1.92%  "
Set-TextValue $ws "D19" "14.53"
$ws.Range("E19").Value = "  +0.61%  "
Set-TextValue $ws "D20" "8.97"
$ws.Range("E20").Value = "  -4.41%  "
Set-TextValue $ws "D21" "383.52"
$ws.Range("E21").Value = "  -1.89%  "
Set-TextValue $ws "D22" "0.571"
$ws.Range("E22").Value = "  +0.88%  "
Set-TextValue $ws "D23" "75.26"
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "3.558.36"
$ws.Range("E25").Value = "  -1.79%  "
$ws.Range("E26").Value = "  -3.22%  "
Set-TextValue $ws "D27" "0.180"
$ws.Range("E27").Value = "  +0.30%  "
Set-TextValue $ws "D28" "7.73"
$ws.Range("E28").Value = "  +1.21%  "
$ws.Range("E29").Value = "  -0.13%  "
Set-TextValue $ws "D30" "7.96"
$ws.Range("E30").Value = "  -3.48%  "
Set-TextValue $ws "D31" "2.13"
$ws.Range("E31").Value = "  -0.90%  "
Set-TextValue $ws "D32" "0.999"
$ws.Range("E32").Value = "  -0.06%  "
Set-TextValue $ws "D33" "1.34"
$ws.Range("E33").Value = "  -1.98%  "
Set-TextValue $ws "D34" "23.25"
$ws.Range("E34").Value = "  -1.70%  "
Set-TextValue $ws "D35" "5.46"
$ws.Range("E35").Value = "  +2.83%  "
Set-TextValue $ws "D36" "1.62"
$ws.Range("E36").Value = "  +1.92%  "
Set-TextValue $ws "D37" "6.93"
$ws.Range("E37").Value = "  -2.99%  "
Set-TextValue $ws "D38" "31.22"
$ws.Range("E38").Value = "  -1.66%  "
Set-TextValue $ws "D39" "168.54"
$ws.Range("E39").Value = "  -0.80%  "
$ws.Range("D40").Value = "3.462.47"
$ws.Range("E40").Value = "  -1.65%  "
Set-TextValue $ws "D41" "0.0783"
$ws.Range("E41").Value = "  +2.00%  "
Set-TextValue $ws "D42" "42.80"
$ws.Range("E42").Value = "  +1.06%  "
Set-TextValue $ws "D43" "0.779"
$ws.Range("E43").Value = "  -2.53%  "
Set-TextValue $ws "D44" "4.41"
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("E45").Value = "  -3.17%  "
Set-TextValue $ws "D46" "1.18"
$ws.Range("E46").Value = "  -3.04%  "
$ws.Range("D47").Value = "2.541.64"
$ws.Range("E47").Value = "  -3.16%  "
Set-TextValue $ws "D48" "6.90"
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("E49").Value = "  -4.65%  "
Set-TextValue $ws "D50" "22.66"
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("E51").Value = "  -0.20%  "
